# This script rewrites Sheet1 from a small "Year / two rows" summary table
# into a wide, unpivoted table with a header row (Unnamed: 0 .. Unnamed: 2,
# 2019, 2018) and four data rows underneath, matching the data-wrangling
# change described in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Wipe out everything that is not part of the new layout (old A2:A5
#    labels, and the old B/C/D row2+row4 values).
# ---------------------------------------------------------------------
$ws.Range("A2:A5").Clear()
$ws.Range("B2:D2").Clear()
$ws.Range("B4:D4").Clear()

# ---------------------------------------------------------------------
# 2. Build a throw-away "text" seed cell so we can clone a text-typed
#    number format onto every cell that must store a numeric-looking
#    value ("2019", "362", "2,082", ...) as literal text rather than a
#    parsed number, exactly like the original authoring tool produced.
#    This must happen BEFORE the header style is spread across the row
#    so that the header cells (B1/D1) end up with both the bold/border
#    look AND text-typed contents.
# ---------------------------------------------------------------------
$seed = $ws.Range("Z100")
$seed.NumberFormat = "@"
$seed.Value = "x"
$seed.Copy()

$numericLookingCells = @("B1", "D1", "B2", "D2", "B3", "D3", "B4", "D4", "B5", "D5")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------
# 3. Fill in the header row values.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("B1").Value = "2019"
$ws.Range("C1").Value = "Unnamed: 1"
$ws.Range("D1").Value = "2018"
$ws.Range("E1").Value = "Unnamed: 2"

# ---------------------------------------------------------------------
# 4. Fill in the four data rows (columns B and D only).
# ---------------------------------------------------------------------
$ws.Range("B2").Value = "362"
$ws.Range("D2").Value = "2,082"

$ws.Range("B3").Value = "494"
$ws.Range("D3").Value = "351"

$ws.Range("B4").Value = "169"
$ws.Range("D4").Value = "447"

$ws.Range("B5").Value = "1,153"
$ws.Range("D5").Value = "3,920"

# ---------------------------------------------------------------------
# 5. Spread the existing header style (bold, centered, bordered) that
#    lives on A1 across the rest of the header row (B1:E1). The text
#    values already stored in B1/D1 stay text-typed even though the
#    number format that formatting carries is "General".
# ---------------------------------------------------------------------
$ws.Range("A1").Copy()
foreach ($addr in @("B1", "C1", "D1", "E1")) {
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------
# 6. Restore the plain/default style on the non-header text cells.
# ---------------------------------------------------------------------
foreach ($addr in @("B2", "D2", "B3", "D3", "B4", "D4", "B5", "D5")) {
    $ws.Range($addr).Style = "Normal"
}

# ---------------------------------------------------------------------
# 7. Clean up the temporary seed cell so it leaves no trace in the
#    final worksheet.
# ---------------------------------------------------------------------
$seed.Clear()

Write-Host "Worksheet rewritten."
